$wb = $excel.ActiveWorkbook

$wsAdmin = $wb.Worksheets.Item("Admin")
$wsJira = $wb.Worksheets.Item("Jira")

# Admin sheet: D2 value changed
$wsAdmin.Range("D2").Value = "1342182A"

# Jira sheet: B2, A3, B3 values changed
$wsJira.Range("B2").Value = "C740395689588328E5DA6BCCD7E88F35"
$wsJira.Range("A3").Value = "PersonalDetails"

# B3's new text ("10600") looks numeric, so force it to stay text
# (leading apostrophe = typed-as-text) and reset style so no extra
# number-format style gets introduced.
$wsJira.Range("B3").Value = "'10600"
$wsJira.Range("B3").Style = "Normal"
